# Austria_JaehrlicheSterbetafeln_1947-2022.xlsx fix:
#  - rename the "row.names" header/table column to "Alter" on all three sheets
#  - convert the age index column (A5:A105, values "1".."101" as text) to the
#    actual numeric age (0..100) on all three sheets

$wb = $excel.ActiveWorkbook

for ($s = 1; $s -le $wb.Worksheets.Count(); $s++) {
    $ws = $wb.Worksheets.Item($s)

    # Header: "row.names" -> "Alter" (also updates the table's column name)
    $ws.Range("A4").Value = "Alter"

    # Age column: replace the text labels with the real numeric age (0-based)
    for ($r = 5; $r -le 105; $r++) {
        $age = $r - 5
        $ws.Cells.Item($r, 1).Value = $age
    }
}
